# Updates cryptos list values (price + 1h volume change) per GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.785.29"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "2.807.91"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.563"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.50%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "3.249.31"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").Value = "2.802.80"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.915"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.29%  "
$ws.Range("D18").Value = "51.627.61"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.59%  "
$ws.Range("E20").Value = "  -4.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").Value = "0.0₃0987"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0495"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +20.51%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("E35").Value = "  +9.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0843"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.116"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "126.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.45%  "
$ws.Range("D46").Value = "2.085.77"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("E49").Value = "  +6.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.970"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.04%  "
$ws.Range("E51").Value = "  +1.31%  "
